{"js": "// \"corrected a few typos\" \u2014 design-discussion.docx\n//\n// 1. \"-Having Discount as a class\" -> \"- Having Discount as a class\"\n//    (missing space after the leading dash)\n// 2. Remove the stray duplicate bullet\n//    \"- Use proper naming convention (cardExpiration instead of CardExpiration)\"\n//    (its content already appears verbatim a few lines below, under Design 2's\n//    Cons, as \"- Variable names could be better\" context \u2014 here it was simply\n//    a leftover duplicate bullet under Design 1's Cons and gets deleted\n//    entirely, paragraph and all).\n// 3. \"...in all of our diagrams credit cards have...\" -> \"...diagrams, credit cards have...\"\n//    (missing comma)\n// 4. \"...our final design \u2013  Rewards...\" -> \"...our final design \u2013 Rewards...\"\n//    (double space after the dash collapsed to one)\n// 5. \"Lessons learnt in teamwork:\" -> \"Lessons learnt in team work:\"\n//    (split into two words)\n\nconst body = context.document.body;\n\n// --- Edit 1: add the missing space after the leading dash ---------------\nconst r1 = body.search(\"-Having Discount as a class\", { matchCase: true });\nr1.load(\"items\");\nawait context.sync();\nif (r1.items.length !== 1) {\n  throw new Error(`Edit 1: expected 1 match, found ${r1.items.length}`);\n}\nr1.items[0].insertText(\"- Having Discount as a class\", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Edit 2: delete the duplicate \"cardExpiration\" bullet paragraph -----\nconst r2 = body.search(\n  \"- Use proper naming convention (cardExpiration instead of CardExpiration)\",\n  { matchCase: true }\n);\nr2.load(\"items,paragraphs\");\nawait context.sync();\nif (r2.items.length !== 1) {\n  throw new Error(`Edit 2: expected 1 match, found ${r2.items.length}`);\n}\nconst paraToDelete = r2.items[0].paragraphs.getFirst();\nparaToDelete.delete();\nawait context.sync();\n\n// --- Edit 3: insert the missing comma ------------------------------------\nconst r3 = body.search(\"in all of our diagrams credit cards have\", { matchCase: true });\nr3.load(\"items\");\nawait context.sync();\nif (r3.items.length !== 1) {\n  throw new Error(`Edit 3: expected 1 match, found ${r3.items.length}`);\n}\nr3.items[0].insertText(\"in all of our diagrams, credit cards have\", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Edit 4: collapse the double space after the dash --------------------\nconst r4 = body.search(\"design \\u2013  Rewards\", { matchCase: true });\nr4.load(\"items\");\nawait context.sync();\nif (r4.items.length !== 1) {\n  throw new Error(`Edit 4: expected 1 match, found ${r4.items.length}`);\n}\nr4.items[0].insertText(\"design \\u2013 Rewards\", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Edit 5: split \"teamwork\" into \"team work\" ---------------------------\nconst r5 = body.search(\"Lessons learnt in teamwork\", { matchCase: true });\nr5.load(\"items\");\nawait context.sync();\nif (r5.items.length !== 1) {\n  throw new Error(`Edit 5: expected 1 match, found ${r5.items.length}`);\n}\nr5.items[0].insertText(\"Lessons learnt in team work\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# \"corrected a few typos\" \u2014 design-discussion.docx\n#\n# 1. \"-Having Discount as a class\" -> \"- Having Discount as a class\"\n#    (missing space after the leading dash)\n# 2. Remove the stray duplicate bullet paragraph\n#    \"- Use proper naming convention (cardExpiration instead of CardExpiration)\"\n# 3. \"...in all of our diagrams credit cards have...\" -> \"...diagrams, credit cards have...\"\n#    (missing comma)\n# 4. \"...our final design -  Rewards...\" -> \"...our final design - Rewards...\"\n#    (double space after the dash collapsed to one)\n# 5. \"Lessons learnt in teamwork:\" -> \"Lessons learnt in team work:\"\n#    (split into two words)\n\n$d = $word.ActiveDocument\n\nfunction Replace-OnceOrFail($FindText, $ReplaceText, $Label) {\n    $range = $d.Content\n    $ok = $range.Find.Execute($FindText, $false, $false, $false, $false, $false, $true, 1, $false, $ReplaceText, 2)\n    if (-not $ok) {\n        throw \"$Label`: text not found -> $FindText\"\n    }\n}\n\n# --- Edit 1: add the missing space after the leading dash -------------------\nReplace-OnceOrFail \"-Having Discount as a class\" \"- Having Discount as a class\" \"Edit 1\"\n\n# --- Edit 2: delete the duplicate \"cardExpiration\" bullet paragraph ---------\n$target = \"- Use proper naming convention (cardExpiration instead of CardExpiration)\"\n$deleted = $false\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $target) {\n        $p.Range.Delete()\n        $deleted = $true\n        break\n    }\n}\nif (-not $deleted) {\n    throw \"Edit 2: paragraph not found -> $target\"\n}\n\n# --- Edit 3: insert the missing comma ---------------------------------------\nReplace-OnceOrFail \"in all of our diagrams credit cards have\" \"in all of our diagrams, credit cards have\" \"Edit 3\"\n\n# --- Edit 4: collapse the double space after the dash -----------------------\n$dash = [char]0x2013\n$edit4Find = \"design \" + $dash + \"  Rewards\"\n$edit4Replace = \"design \" + $dash + \" Rewards\"\nReplace-OnceOrFail $edit4Find $edit4Replace \"Edit 4\"\n\n# --- Edit 5: split \"teamwork\" into \"team work\" ------------------------------\nReplace-OnceOrFail \"Lessons learnt in teamwork\" \"Lessons learnt in team work\" \"Edit 5\"\n"}
